# Auto-generated edit script: updates cryptos price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.739.83"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.628.64"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'214.51"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.257"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("E10").Value = "  -4.81%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "1.631.42"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "1.852.57"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Value = "0.0₃0761"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "'62.69"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "25.742.73"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "'193.83"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'9.92"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "'139.63"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "'15.46"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "'0.0488"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'3.24"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'0.546"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "1.107.92"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "'99.94"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "'0.800"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "1.758.86"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "'54.92"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'7.74"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.417"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0501"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.36"
$ws.Range("E51").Value = "  +2.44%  "
